# Insert a new price-report row for "Ají" (Vega Monumental Concepción) as row 147,
# pushing the existing rows 147-236 down to 148-237.
#
# Net effect on the sheet: dimension A1:R236 -> A1:R237, and a brand-new
# observation is recorded at row 147 (date 2023-10-17 / serial 45216,
# volume 80, price range 30000-32000, weighted avg 30750, "$/caja 10 kilos",
# 3075 $/Kg over 10 Kg) while every previously existing row keeps its data,
# just shifted one row further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 147; this shifts rows 147:236
# down to 148:237 and keeps all of their values/formatting intact.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new observation.
$ws.Range("A147").Value = 11
$ws.Range("B147").Value = "Vega Monumental Concepción"
$ws.Range("C147").Value = "Bíobío"
$ws.Range("D147").Value = 45216
$ws.Range("E147").Value = 8
$ws.Range("F147").Value = 100112021
$ws.Range("G147").Value = "Ají"
$ws.Range("H147").Value = "Inferno"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 80
$ws.Range("K147").Value = 30000
$ws.Range("L147").Value = 32000
$ws.Range("M147").Value = 30750
$ws.Range("N147").Value = "$/caja 10 kilos"
$ws.Range("O147").Value = "Región de Arica y Parinacota"
$ws.Range("P147").Value = 3075
$ws.Range("Q147").Value = 10
$ws.Range("R147").Value = "Hortaliza"
